# Weekly crime data update for cs-en-us-120pct.xlsx
# - Bumps the report volume/number and the covered week's date range.
# - Refreshes the weekly/28-day/YTD/2-year crime-count + %-change figures
#   for the period 6/16/2025 - 6/22/2025 (previously 6/9/2025 - 6/15/2025).
#
# Many of the numeric cells in this report render as literal placeholder
# text ("0" or "***.*") instead of a real number/percentage when a count is
# zero or a percent-change is undefined (division by zero) in the source
# report. Set-Cell below mirrors that: for string targets it forces the
# cell to Text format first so Excel does not auto-coerce a numeric-looking
# string like "0" back into a number, then restores the cell's real
# (General/number) format afterwards so the stored format matches the
# report's normal styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell {
    param($sheet, $addr, $value, $fmt)
    $rng = $sheet.Range($addr)
    if ($value -is [string]) {
        # Pre-format as Text so numeric-looking strings (e.g. "0") are not
        # auto-converted to a number by the COM Value setter.
        $rng.NumberFormat = "@"
        $rng.Value = $value
        $rng.NumberFormat = $fmt
    } else {
        $rng.Value = $value
        $rng.NumberFormat = $fmt
    }
}

Set-Cell $ws "A8" 'Volume 32   Number  25' 'General'
Set-Cell $ws "C9" 'Report Covering the Week  6/16/2025  Through  6/22/2025' 'General'
Set-Cell $ws "D14" 1 '#,##0'
Set-Cell $ws "E14" -100 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "G14" 1 '#,##0'
Set-Cell $ws "H14" -100 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "J14" 1 '#,##0'
Set-Cell $ws "K14" 100 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "N14" -81.818181818181 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "G15" '0' 'General'
Set-Cell $ws "H15" '***.*' 'General'
Set-Cell $ws "L15" 30 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "N15" -55.172413793103 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "C16" 1 '#,##0'
Set-Cell $ws "E16" -50 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "F16" 9 '#,##0'
Set-Cell $ws "G16" 8 '#,##0'
Set-Cell $ws "H16" 12.5 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "I16" 66 '#,##0'
Set-Cell $ws "J16" 66 '#,##0'
Set-Cell $ws "K16" 0 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "L16" -19.512195121951 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "M16" -48.031496062992 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "N16" -85.714285714285 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "C17" 8 '#,##0'
Set-Cell $ws "D17" 9 '#,##0'
Set-Cell $ws "E17" -11.111111111111 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "F17" 36 '#,##0'
Set-Cell $ws "H17" -7.692307692307 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "I17" 204 '#,##0'
Set-Cell $ws "J17" 235 '#,##0'
Set-Cell $ws "K17" -13.191489361702 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "L17" 4.615384615384 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "M17" 34.210526315789 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "N17" -43.490304709141 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "C18" 1 '#,##0'
Set-Cell $ws "D18" '0' 'General'
Set-Cell $ws "E18" '***.*' 'General'
Set-Cell $ws "F18" 3 '#,##0'
Set-Cell $ws "G18" 8 '#,##0'
Set-Cell $ws "H18" -62.5 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "I18" 49 '#,##0'
Set-Cell $ws "K18" 0 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "L18" -33.783783783783 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "M18" -57.391304347826 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "N18" -94.067796610169 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "C19" 11 '#,##0'
Set-Cell $ws "D19" 12 '#,##0'
Set-Cell $ws "E19" -8.333333333333 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "F19" 29 '#,##0'
Set-Cell $ws "G19" 45 '#,##0'
Set-Cell $ws "H19" -35.555555555555 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "I19" 142 '#,##0'
Set-Cell $ws "J19" 206 '#,##0'
Set-Cell $ws "K19" -31.067961165048 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "L19" -24.064171122994 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "M19" -24.867724867724 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "N19" -44.961240310077 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "C20" 1 '#,##0'
Set-Cell $ws "D20" 6 '#,##0'
Set-Cell $ws "E20" -83.333333333333 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "F20" 5 '#,##0'
Set-Cell $ws "G20" 18 '#,##0'
Set-Cell $ws "H20" -72.222222222222 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "I20" 23 '#,##0'
Set-Cell $ws "J20" 58 '#,##0'
Set-Cell $ws "K20" -60.344827586206 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "L20" -61.016949152542 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "M20" -73.563218390804 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "N20" -96.40062597809 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "C21" 22 '#,##0'
Set-Cell $ws "D21" 30 '#,##0'
Set-Cell $ws "E21" -26.666666666666 '#,##0.00;"-"#,##0.00'
Set-Cell $ws "F21" 83 '#,##0'
Set-Cell $ws "G21" 119 '#,##0'
Set-Cell $ws "H21" -30.252100840336 '#,##0.00;"-"#,##0.00'
Set-Cell $ws "I21" 499 '#,##0'
Set-Cell $ws "J21" 626 '#,##0'
Set-Cell $ws "K21" -20.287539936102 '#,##0.00;"-"#,##0.00'
Set-Cell $ws "L21" -18.861788617886 '#,##0.00;"-"#,##0.00'
Set-Cell $ws "M21" -27.470930232558 '#,##0.00;"-"#,##0.00'
Set-Cell $ws "N21" -80.703789636504 '#,##0.00;"-"#,##0.00'
Set-Cell $ws "C23" 1 '#,##0'
Set-Cell $ws "D23" 2 '#,##0'
Set-Cell $ws "E23" -50 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "F23" 5 '#,##0'
Set-Cell $ws "H23" -16.666666666666 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "I23" 29 '#,##0'
Set-Cell $ws "J23" 39 '#,##0'
Set-Cell $ws "K23" -25.641025641025 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "L23" -35.555555555555 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "M23" 52.631578947368 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "C24" 19 '#,##0'
Set-Cell $ws "D24" 22 '#,##0'
Set-Cell $ws "E24" -13.636363636363 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "F24" 94 '#,##0'
Set-Cell $ws "G24" 87 '#,##0'
Set-Cell $ws "H24" 8.045977011494 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "I24" 626 '#,##0'
Set-Cell $ws "J24" 633 '#,##0'
Set-Cell $ws "K24" -1.105845181674 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "L24" 9.440559440559 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "M24" -4.281345565749 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "C25" 12 '#,##0'
Set-Cell $ws "D25" 3 '#,##0'
Set-Cell $ws "E25" 300 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "F25" 47 '#,##0'
Set-Cell $ws "G25" 28 '#,##0'
Set-Cell $ws "H25" 67.857142857142 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "I25" 334 '#,##0'
Set-Cell $ws "J25" 288 '#,##0'
Set-Cell $ws "K25" 15.972222222222 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "L25" 46.491228070175 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "C26" 17 '#,##0'
Set-Cell $ws "D26" 22 '#,##0'
Set-Cell $ws "E26" -22.727272727272 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "F26" 60 '#,##0'
Set-Cell $ws "G26" 74 '#,##0'
Set-Cell $ws "H26" -18.918918918918 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "I26" 350 '#,##0'
Set-Cell $ws "J26" 351 '#,##0'
Set-Cell $ws "K26" -0.2849002849 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "L26" 4.790419161676 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "M26" -38.38028169014 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "C27" '0' 'General'
Set-Cell $ws "G27" '0' 'General'
Set-Cell $ws "H27" '***.*' 'General'
Set-Cell $ws "L27" 66.666666666666 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "C28" 2 '#,##0'
Set-Cell $ws "E28" 0 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "F28" 4 '#,##0'
Set-Cell $ws "G28" 5 '#,##0'
Set-Cell $ws "H28" -20 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "I28" 43 '#,##0'
Set-Cell $ws "J28" 38 '#,##0'
Set-Cell $ws "K28" 13.157894736842 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "L28" 0 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "D29" '0' 'General'
Set-Cell $ws "E29" '***.*' 'General'
Set-Cell $ws "L29" -83.333333333333 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "M29" -88.235294117647 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "N29" -95.918367346938 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "D30" '0' 'General'
Set-Cell $ws "E30" '***.*' 'General'
Set-Cell $ws "L30" -81.818181818181 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "M30" -86.666666666666 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "N30" -94.871794871794 '#,##0.0;"-"#,##0.0'
Set-Cell $ws "D31" '0' 'General'
Set-Cell $ws "E31" '***.*' 'General'
